# Renal Tubulopathies workbook: add a "metadata" tab and refresh the
# "data" sheet's per-row query timestamps (column F) to the later
# panel-fetch run.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update the time_taken column (F2:F36) on the existing "data" sheet
#    with the refreshed query timestamps.
# ---------------------------------------------------------------------
$timestamps = @(
    "2021-10-05 14:35:36.503884",
    "2021-10-05 14:35:36.503892",
    "2021-10-05 14:35:36.503896",
    "2021-10-05 14:35:36.503899",
    "2021-10-05 14:35:36.503901",
    "2021-10-05 14:35:36.503904",
    "2021-10-05 14:35:36.503907",
    "2021-10-05 14:35:36.503909",
    "2021-10-05 14:35:36.503912",
    "2021-10-05 14:35:36.503915",
    "2021-10-05 14:35:36.503917",
    "2021-10-05 14:35:36.503920",
    "2021-10-05 14:35:36.503922",
    "2021-10-05 14:35:36.503925",
    "2021-10-05 14:35:36.503927",
    "2021-10-05 14:35:36.503930",
    "2021-10-05 14:35:36.503933",
    "2021-10-05 14:35:36.503935",
    "2021-10-05 14:35:36.503938",
    "2021-10-05 14:35:36.503941",
    "2021-10-05 14:35:36.503943",
    "2021-10-05 14:35:36.503946",
    "2021-10-05 14:35:36.503948",
    "2021-10-05 14:35:36.503951",
    "2021-10-05 14:35:36.503954",
    "2021-10-05 14:35:36.503956",
    "2021-10-05 14:35:36.503959",
    "2021-10-05 14:35:36.503961",
    "2021-10-05 14:35:36.503964",
    "2021-10-05 14:35:36.503966",
    "2021-10-05 14:35:36.503969",
    "2021-10-05 14:35:36.503971",
    "2021-10-05 14:35:36.503974",
    "2021-10-05 14:35:36.503977",
    "2021-10-05 14:35:36.503980"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# ---------------------------------------------------------------------
# 2. Add a new "metadata" worksheet after "data" describing the panel
#    query that produced this export.
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Renal Tubulopathies"
$metaSheet.Range("C2").Value = 200
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.44"
$metaSheet.Range("E2").Value = "2021-09-16T21:12:23.131806Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:36.500035"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/200/?format=json"

# Match the header/index-column styling used on the "data" sheet (bold,
# bordered, centered) by copying formats from the equivalent cells.
$dataSheet.Range("B1:F1").Copy() | Out-Null
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null

$dataSheet.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
